# Update crypto price/volume data in Sheet1 per the Oct 16 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '27.888.24'
$ws.Range('E2').Value = '  +2.87%  '
$ws.Range('D3').Value = '1.569.92'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').Value = "'0.989"
$ws.Range('E4').Value = '  -1.94%  '
$ws.Range('D5').Value = "'211.47"
$ws.Range('E5').Value = '  +0.32%  '
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('E7').Value = '  -1.92%  '
$ws.Range('D8').Value = "'23.20"
$ws.Range('E8').Value = '  +5.76%  '
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('E11').Value = '  +1.51%  '
$ws.Range('D12').Value = '1.795.18'
$ws.Range('E12').Value = '  +0.29%  '
$ws.Range('D13').Value = '1.570.67'
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('E14').Value = '  -0.61%  '
$ws.Range('D15').Value = "'0.519"
$ws.Range('E15').Value = '  +0.09%  '
$ws.Range('D16').Value = '27.853.92'
$ws.Range('D17').Value = "'63.44"
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('D18').Value = "'230.57"
$ws.Range('E18').Value = '  +7.38%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = "'7.44"
$ws.Range('E21').Value = '  -1.95%  '
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = "'1.92"
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').Value = "'151.13"
$ws.Range('E25').Value = '  -1.93%  '
$ws.Range('D26').Value = "'15.23"
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E33').Value = '  -1.86%  '
$ws.Range('D34').Value = '1.419.06'
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('E36').Value = '  -4.43%  '
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').Value = "'0.542"
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').Value = "'2.40"
$ws.Range('E40').Value = '  +1.40%  '
$ws.Range('E41').Value = '  -0.22%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'0.989"
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'5.61"
$ws.Range('E43').Value = '  -3.68%  '
$ws.Range('E44').Value = '  +5.54%  '
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('D46').Value = "'63.89"
$ws.Range('E46').Value = '  -0.99%  '
$ws.Range('D47').Value = '1.705.18'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('D48').Value = "'86.63"
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').Value = '0.0₇0984'
$ws.Range('E50').Value = '  -3.92%  '
$ws.Range('D51').Value = "'39.65"
$ws.Range('E51').Value = '  +17.09%  '
